$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (24 -> 22 rows total)
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(23).Delete()

# Row 2
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "Cappucu1no Pq 150 m1 1 8"
$ws.Range("B2").Value = 500
$ws.Range("C2").Value = 8
$ws.Range("D2").Value = 500
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "otros"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "11/01/2024"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "por definir"

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "Hatsu"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 49900
$ws.Range("D3").Value = 49900
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "bebida sin alcohol"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "desconocida"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "Pricesmart"

# Row 4
$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "Filete Pech"
$ws.Range("B4").Value = 294588
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 58
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "otros"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "desconocida"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "Pricesmart"

# Row 5
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "Lomo Cerdo"
$ws.Range("B5").Value = 86320
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 48
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "carnes rojas"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "desconocida"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "Pricesmart"

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "Ques!to"
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 22
$ws.Range("D6").Value = 22
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "otros"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "desconocida"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "Pricesmart"

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "HuevosAA 60"
$ws.Range("B7").Value = 1
$ws.Range("C7").Value = 32
$ws.Range("D7").Value = 32
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "carnes blancas"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "desconocida"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "Pricesmart"

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = "Fideos"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = 10
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "cereales"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "desconocida"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Pricesmart"

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "Colageno"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = 78
$ws.Range("D9").Value = 78
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "otros"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "desconocida"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Pricesmart"

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "Miel"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = 46
$ws.Range("D10").Value = 46
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "otros"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "desconocida"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "Pricesmart"

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "Mantequilla"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = 52
$ws.Range("D11").Value = 52
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "lácteos"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "desconocida"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Pricesmart"

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "AguaSabor1x"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = 53
$ws.Range("D12").Value = 53
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "bebida sin alcohol"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "desconocida"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "Pricesmart"

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "KSR6dCUPZ4D"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = 82
$ws.Range("D13").Value = 82
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "otros"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "desconocida"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "Pricesmart"

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "Hatsu"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = 49900
$ws.Range("D14").Value = 49900
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "bebida sin alcohol"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "desconocida"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "Pricesmart"

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "Filete Pech"
$ws.Range("B15").Value = 294588
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 58
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "carnes blancas"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "desconocida"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "Pricesmart"

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "Lomo Cerdo"
$ws.Range("B16").Value = 86320
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 48
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "carnes rojas"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "desconocida"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "Pricesmart"

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "Quesito"
$ws.Range("B17").Value = 1
$ws.Range("C17").Value = 22
$ws.Range("D17").Value = 22
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "otros"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "desconocida"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "Pricesmart"

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "HuevosAA 60"
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 32
$ws.Range("D18").Value = 32
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "otros"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "desconocida"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "Pricesmart"

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "Fideos"
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 10
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "otros"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "desconocida"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "Pricesmart"

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "i 333342 Almendra21b"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 46
$ws.Range("D20").Value = 46
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "otros"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "desconocida"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "Pricesmart"

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = "cerveza"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = 82000
$ws.Range("D21").Value = 82000
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "otros"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "desconocida"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "Pricesmart"

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "Cappucu1no Pq 150 m1 1 8 500 8"
$ws.Range("B22").Value = 1
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "500"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "500"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "otros"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "11/01/2024"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "por definir"
